$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1917808219178082
$ws.Range("C2").Value = 0.564134495641345
$ws.Range("J2").Value = 0.009962640099626401
$ws.Range("P2").Value = 0.1432129514321295
$ws.Range("S2").Value = 0.09090909090909091
$ws.Range("B3").Value = 0.008403361344537815
$ws.Range("C3").Value = 0.04831932773109244
$ws.Range("J3").Value = 0.01470588235294118
$ws.Range("P3").Value = 0.7478991596638656
$ws.Range("S3").Value = 0.180672268907563
$ws.Range("J4").Value = 0.06024096385542169
$ws.Range("P4").Value = 0.6265060240963856
$ws.Range("S4").Value = 0.3132530120481928
$ws.Range("B6").Value = 0.05573248407643312
$ws.Range("D6").Value = 0.009554140127388535
$ws.Range("F6").Value = 0.07165605095541401
$ws.Range("J6").Value = 0.2579617834394904
$ws.Range("O6").Value = 0.01592356687898089
$ws.Range("Q6").Value = 0.1624203821656051
$ws.Range("R6").Value = 0.05732484076433121
$ws.Range("S6").Value = 0.3694267515923567
$ws.Range("B7").Value = 0.09263157894736843
$ws.Range("D7").Value = 0.01894736842105263
$ws.Range("E7").Value = 0.002105263157894737
$ws.Range("F7").Value = 0.04210526315789474
$ws.Range("J7").Value = 0.12
$ws.Range("O7").Value = 0.01894736842105263
$ws.Range("Q7").Value = 0.1726315789473684
$ws.Range("R7").Value = 0.1031578947368421
$ws.Range("S7").Value = 0.4294736842105263
$ws.Range("B8").Value = 0.09923664122137404
$ws.Range("D8").Value = 0.01221374045801527
$ws.Range("E8").Value = 0.002290076335877863
$ws.Range("F8").Value = 0.06183206106870229
$ws.Range("J8").Value = 0.1022900763358779
$ws.Range("O8").Value = 0.01374045801526718
$ws.Range("Q8").Value = 0.1702290076335878
$ws.Range("R8").Value = 0.09770992366412214
$ws.Range("S8").Value = 0.4404580152671756
$ws.Range("B9").Value = 0.1030927835051546
$ws.Range("D9").Value = 0.01202749140893471
$ws.Range("E9").Value = 0.001718213058419244
$ws.Range("F9").Value = 0.07903780068728522
$ws.Range("J9").Value = 0.09793814432989691
$ws.Range("O9").Value = 0.01890034364261168
$ws.Range("Q9").Value = 0.1941580756013746
$ws.Range("R9").Value = 0.07903780068728522
$ws.Range("S9").Value = 0.4140893470790378
$ws.Range("B10").Value = 0.11121103117506
$ws.Range("D10").Value = 0.01438848920863309
$ws.Range("E10").Value = 0.0005995203836930455
$ws.Range("F10").Value = 0.07254196642685851
$ws.Range("J10").Value = 0.1019184652278178
$ws.Range("O10").Value = 0.01528776978417266
$ws.Range("Q10").Value = 0.2113309352517986
$ws.Range("R10").Value = 0.08992805755395683
$ws.Range("S10").Value = 0.3827937649880096
$ws.Range("G11").Value = 0.1549497847919656
$ws.Range("J11").Value = 0.07317073170731707
$ws.Range("K11").Value = 0.1908177905308465
$ws.Range("L11").Value = 0.5581061692969871
$ws.Range("S11").Value = 0.02295552367288379
$ws.Range("G12").Value = 0.748792270531401
$ws.Range("J12").Value = 0.1739130434782609
$ws.Range("K12").Value = 0.007246376811594203
$ws.Range("L12").Value = 0.03864734299516908
$ws.Range("S12").Value = 0.03140096618357488
$ws.Range("G13").Value = 0.6724137931034483
$ws.Range("J13").Value = 0.2844827586206897
$ws.Range("S13").Value = 0.04310344827586207
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.375
$ws.Range("S14").Value = 0.125
$ws.Range("F15").Value = 0.02711323763955343
$ws.Range("H15").Value = 0.1610845295055821
$ws.Range("I15").Value = 0.08452950558213716
$ws.Range("J15").Value = 0.3763955342902711
$ws.Range("K15").Value = 0.04944178628389154
$ws.Range("M15").Value = 0.0127591706539075
$ws.Range("O15").Value = 0.06379585326953748
$ws.Range("S15").Value = 0.2248803827751196
$ws.Range("F16").Value = 0.03326810176125244
$ws.Range("H16").Value = 0.1761252446183953
$ws.Range("I16").Value = 0.08023483365949119
$ws.Range("J16").Value = 0.4031311154598826
$ws.Range("K16").Value = 0.09001956947162426
$ws.Range("M16").Value = 0.01565557729941291
$ws.Range("N16").Value = 0.003913894324853229
$ws.Range("O16").Value = 0.06066536203522505
$ws.Range("S16").Value = 0.136986301369863
$ws.Range("F17").Value = 0.0217566478646253
$ws.Range("H17").Value = 0.1635777598710717
$ws.Range("I17").Value = 0.09669621273166801
$ws.Range("J17").Value = 0.4174053182917002
$ws.Range("K17").Value = 0.07977437550362611
$ws.Range("M17").Value = 0.016116035455278
$ws.Range("O17").Value = 0.06688154713940371
$ws.Range("S17").Value = 0.1377921031426269
$ws.Range("F18").Value = 0.01964285714285714
$ws.Range("H18").Value = 0.1714285714285714
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.4375
$ws.Range("K18").Value = 0.075
$ws.Range("M18").Value = 0.01607142857142857
$ws.Range("N18").Value = 0.001785714285714286
$ws.Range("O18").Value = 0.06428571428571428
$ws.Range("S18").Value = 0.1142857142857143
$ws.Range("F19").Value = 0.01752433936022253
$ws.Range("H19").Value = 0.2303198887343533
$ws.Range("I19").Value = 0.08845618915159945
$ws.Range("J19").Value = 0.3452016689847009
$ws.Range("K19").Value = 0.09068150208623088
$ws.Range("M19").Value = 0.01974965229485396
$ws.Range("N19").Value = 0.001390820584144645
$ws.Range("O19").Value = 0.0717663421418637
$ws.Range("S19").Value = 0.1349095966620306
